$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "mean_years_of_schooling"
$ws.Range("B5").Value = "average years of sccooling per country"
$ws.Range("C5").Value = "1870-2017"
$ws.Range("D5").Value = "Yearly"

$ws.Range("E5").Select()
